$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh (GitHub Actions run, Sat Apr 29 18:09:44 UTC 2023).
# Price (D) and Volume(1h) (E) columns hold plain text (e.g. "29.436.30",
# "  +0.65%  ") rather than numbers, so force Text format before writing so
# Excel doesn't auto-convert the digit strings to numeric values.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.436.30"
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "1.913.47"
$ws.Range("E3").Value = "  +1.16%  "
$ws.Range("D4").Value = "1.011"
$ws.Range("E4").Value = "  +0.89%  "
$ws.Range("D5").Value = "325.56"
$ws.Range("E5").Value = "  +0.91%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "0.4811"
$ws.Range("E7").Value = "  +1.27%  "
$ws.Range("D8").Value = "0.4077"
$ws.Range("E8").Value = "  +0.85%  "
$ws.Range("D9").Value = "0.08234"
$ws.Range("E9").Value = "  +2.59%  "
$ws.Range("D10").Value = "1.021"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").Value = "23.40"
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").Value = "1.903.31"
$ws.Range("E12").Value = "  +1.23%  "
$ws.Range("D13").Value = "6.043"
$ws.Range("E13").Value = "  +2.14%  "
$ws.Range("D14").Value = "7.223"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "91.06"
$ws.Range("E15").Value = "  +2.16%  "

# Rows 16/17 swapped (TRON now ranks above BinanceUSD)
$ws.Range("B16").Value = "TRON"
$ws.Range("C16").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D16").Value = "0.06814"
$ws.Range("E16").Value = "  +2.78%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  +1.05%  "
$ws.Range("D19").Value = "17.71"
$ws.Range("E19").Value = "  +1.12%  "
$ws.Range("D20").Value = "1.009"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "29.460.77"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "5.630"
$ws.Range("E22").Value = "  +2.43%  "
$ws.Range("D23").Value = "11.77"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "2.200"
$ws.Range("E24").Value = "  +1.22%  "
$ws.Range("D25").Value = "2.157.84"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").Value = "6.613"
$ws.Range("E26").Value = "  +11.39%  "
$ws.Range("D27").Value = "157.13"
$ws.Range("E27").Value = "  +1.67%  "
$ws.Range("D28").Value = "20.09"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "2.108"
$ws.Range("E29").Value = "  +1.38%  "
$ws.Range("D30").Value = "120.25"
$ws.Range("E30").Value = "  +2.44%  "
$ws.Range("D31").Value = "1.023"
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("D32").Value = "0.09561"
$ws.Range("E32").Value = "  +1.44%  "
$ws.Range("D33").Value = "5.572"
$ws.Range("E33").Value = "  +4.41%  "
$ws.Range("D34").Value = "3.556"
$ws.Range("E34").Value = "  +0.84%  "
$ws.Range("D35").Value = "1.367"
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("D36").Value = "0.02280"
$ws.Range("E36").Value = "  +1.42%  "
$ws.Range("D37").Value = "0.06122"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "1.182"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("D39").Value = "8.065"
$ws.Range("E39").Value = "  +1.82%  "

# Rows 40/41 swapped (Aptos now ranks above TheSandbox)
$ws.Range("B40").Value = "Aptos"
$ws.Range("C40").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D40").Value = "10.85"
$ws.Range("E40").Value = "  +8.25%  "
$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").Value = "0.5966"
$ws.Range("E41").Value = "  +2.38%  "

$ws.Range("D42").Value = "0.1853"
$ws.Range("E42").Value = "  +1.17%  "
$ws.Range("D43").Value = "2.411"
$ws.Range("E43").Value = "  +2.01%  "
$ws.Range("D44").Value = "1.281"
$ws.Range("E44").Value = "  -0.58%  "

# Rows 45/46 swapped (EnergySwap now ranks above Cronos)
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "12.46"
$ws.Range("E45").Value = "  +2.49%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "0.07594"
$ws.Range("E46").Value = "  -1.36%  "

$ws.Range("D47").Value = "0.5571"
$ws.Range("E47").Value = "  +1.74%  "
$ws.Range("D48").Value = "1.951"
$ws.Range("E48").Value = "  +2.30%  "
$ws.Range("D49").Value = "117.32"
$ws.Range("E49").Value = "  +3.81%  "
$ws.Range("D50").Value = "2.424"
$ws.Range("E50").Value = "  +4.02%  "
$ws.Range("D51").Value = "72.42"

# Restore default (unstyled) look now that the values are committed as text.
$dataRange.Style = "Normal"
